# Removed Type from imports
# Delete the "Commitment Type" column (column G) from the FundFormula sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FundFormula")

# Select column G (matches the selection seen after the edit: G1:G1048576)
$ws.Range("G1:G1048576").Select() | Out-Null

# Delete the entire column, shifting H:K left to G:J
$ws.Columns.Item(7).Delete() | Out-Null
